# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps on the
# zh-cn and de-de worksheets to reflect a fresh handback run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-20 10:52:48"
$wsZhCn.Range("H2").Value = "2016-03-20 10:53:08"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-20 10:52:52"
$wsDeDe.Range("H2").Value = "2016-03-20 10:53:13"
